$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update terminology from '外籍' (foreign) to '境外' (overseas) in header cells
$ws.Range("C1").Value = "境外學者姓名"
$ws.Range("D1").Value = "境外學者身分（教授、副教授、助理教授或博士後研究員）"

# Update the active selection to D11 as recorded in the workbook view
$ws.Range("D11").Select()
